# Apply updated values to column F ("dSF") on Sheet1, reflecting a
# repull / recalculation of the data (commit: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value  = -5
$ws.Range("F4").Value  = -3
$ws.Range("F6").Value  = -3
$ws.Range("F7").Value  = -5
$ws.Range("F8").Value  = -2
$ws.Range("F9").Value  = -5
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -1
$ws.Range("F17").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = 3
$ws.Range("F36").Value = 2
$ws.Range("F41").Value = 1
$ws.Range("F44").Value = 1
$ws.Range("F47").Value = -6
$ws.Range("F48").Value = -9
